$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from I1 into J1, then set header text.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Route"

# Fill the new "Route" column: every Hoshido-aligned class is "Birthright",
# every Nohr-aligned class is "Conquest" (Nohr Prince(ss) and Nohr Noble play
# into Conquest even though they're listed with the Hoshido block).
$birthrightRows = @(3, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32)
$conquestRows = @(2, 4, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65)

foreach ($row in $birthrightRows) {
    $ws.Cells.Item($row, 10).Value = "Birthright"
}
foreach ($row in $conquestRows) {
    $ws.Cells.Item($row, 10).Value = "Conquest"
}

# Match the saved selection / scroll position from the edited workbook.
[void]$ws.Range("J14").Select()
